$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(227).Insert()

$ws.Cells.Item(227, 1).Value = 4
$ws.Cells.Item(227, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(227, 3).Value = "Los Lagos"
$ws.Cells.Item(227, 4).Value = 44508
$ws.Cells.Item(227, 5).Value = 10
$ws.Cells.Item(227, 6).Value = 100112006
$ws.Cells.Item(227, 7).Value = "Repollo"
$ws.Cells.Item(227, 8).Value = "Crespo record"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 500
$ws.Cells.Item(227, 11).Value = 1200
$ws.Cells.Item(227, 12).Value = 1200
$ws.Cells.Item(227, 13).Value = 1200
$ws.Cells.Item(227, 14).Value = "`$/unidad"
$ws.Cells.Item(227, 15).Value = "Región Metropolitana"
$ws.Cells.Item(227, 16).Value = 1200
$ws.Cells.Item(227, 17).Value = 1
$ws.Cells.Item(227, 18).Value = "Hortaliza"
